$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13234
$ws1.Range("F5").Value = 66
$ws1.Range("F6").Value = 112
$ws1.Range("F8").Value = 60
$ws1.Range("F11").Value = 13177
$ws1.Range("F12").Value = 321
$ws1.Range("F14").Value = 8831
$ws1.Range("F15").Value = 7908
$ws1.Range("F27").Value = 85
$ws1.Range("F28").Value = 351

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13234
$ws4.Range("F6").Value = 66
$ws4.Range("F7").Value = 112
$ws4.Range("F9").Value = 60
$ws4.Range("F12").Value = 13177
$ws4.Range("F13").Value = 321
$ws4.Range("F15").Value = 8831
$ws4.Range("F16").Value = 7908
$ws4.Range("F30").Value = 85
$ws4.Range("F31").Value = 351
